# "added 4wk low sales check"
# Updates the Inventory Coverage (H) and Seasonality Index (L) figures on
# the "Forecast Comparison" sheet to reflect the new 4-week-low sales
# check, and refreshes the derived "Total Forecast (8 Weeks)" figure on
# the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$fc = $wb.Worksheets.Item("Forecast Comparison")
$sm = $wb.Worksheets.Item("Summary")

# Inventory Coverage (column H)
$fc.Range("H2").Value = 5
$fc.Range("H3").Value = 4
$fc.Range("H4").Value = 3
$fc.Range("H5").Value = 2
$fc.Range("H6").Value = 1
$fc.Range("H7").Value = 0

# Seasonality Index (column L)
$fc.Range("L2").Value = 1.2
$fc.Range("L3").Value = 1.01
$fc.Range("L4").Value = 0.99
$fc.Range("L5").Value = 1.04
$fc.Range("L6").Value = 1.12
$fc.Range("L7").Value = 0.83
$fc.Range("L8").Value = 0.84
$fc.Range("L9").Value = 1
$fc.Range("L10").Value = 1
$fc.Range("L11").Value = 0.82
$fc.Range("L12").Value = 1.14
$fc.Range("L13").Value = 1.17
$fc.Range("L14").Value = 0.95
$fc.Range("L15").Value = 0.85
$fc.Range("L16").Value = 0.96
$fc.Range("L17").Value = 0.98

# Summary sheet: "Total Forecast (8 Weeks)" is stored as text, not a
# number, in this workbook. A plain .Value assignment of a numeric-looking
# string gets auto-converted to a number by Excel, so force it to stay
# text (leading apostrophe), then strip the resulting "Text" formatting
# back to Normal so only the value itself changes.
$sm.Range("B10").Formula = "'2"
$sm.Range("B10").NumberFormat = "General"
$sm.Range("B10").Style = "Normal"
